$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format from the last existing data row (A7) so the new
# rows 8 and 9 keep the same look (bold, centered, bordered) without
# introducing any new cell-style entries.
$fmtSource = $ws.Range("A7")

# Row 8: 14/02/2022
$ws.Range("A8").Value = "14/02/2022"
$fmtSource.Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = 3809.882294
$ws.Range("C8").Value = 11358.735837
$ws.Range("D8").Value = 0.54
$ws.Range("E8").Value = 10.65

# Row 9: 15/02/2022
$ws.Range("A9").Value = "15/02/2022"
$fmtSource.Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").Value = 3810.859521
$ws.Range("C9").Value = 11363.298414
$ws.Range("D9").Value = 0.54
$ws.Range("E9").Value = 10.65

$excel.CutCopyMode = $false
